# DEV 13 - Bug Fix and MVC Project, Application, Enquiry
#
# ProjectList.xlsx updates:
#   - Manager column (L) changed from NRIC numbers to the manager's name "Michael"
#     for every project row.
#   - "Acacia Breeze" (row 2) visibility flag fixed from "Hidden" to "Visible".
#   - Project IDs renumbered sequentially (1,2,3,4,5) after the fix-ups.
#   - Two new projects appended: "Savannah Condopark" (Simei) and
#     "Archipelago" (Bedok).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Acacia Breeze ---
$ws.Range("L2").Value = "Michael"
$ws.Range("O2").Value = "Visible"

# --- Row 3: Melville Park -- Project ID renumbered 3 -> 2, manager fixed ---
$ws.Range("A3").Value = 2
$ws.Range("L3").Value = "Michael"

# --- Row 4: Pasir Ris Park -- Project ID renumbered 4 -> 3, manager fixed ---
$ws.Range("A4").Value = 3
$ws.Range("L4").Value = "Michael"

# --- Row 5 (new): Savannah Condopark ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Savannah Condopark"
$ws.Range("C5").Value = "Simei"
$ws.Range("D5").Value = "2-ROOM"
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 100000
$ws.Range("G5").Value = "3-ROOM"
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 200000
$ws.Range("J5").Value = 45775
$ws.Range("K5").Value = 45805
$ws.Range("L5").Value = "Michael"
$ws.Range("M5").Value = 10
$ws.Range("O5").Value = "Visible"

# --- Row 6 (new): Archipelago ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Archipelago"
$ws.Range("C6").Value = "Bedok"
$ws.Range("D6").Value = "2-ROOM"
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 100
$ws.Range("G6").Value = "3-ROOM"
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 45775
$ws.Range("K6").Value = 45805
$ws.Range("L6").Value = "Michael"
$ws.Range("M6").Value = 10
$ws.Range("O6").Value = "Visible"
